# Rename output variable sheets + flip signs of the affected component columns.

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets (dot -> underscore, plus "time" -> "time_Phi") ---
$wb.Worksheets.Item("A.hat").Name = "A_hat"
$wb.Worksheets.Item("B.hat").Name = "B_hat"
$wb.Worksheets.Item("Phi.hat").Name = "Phi_hat"
$wb.Worksheets.Item("time").Name = "time_Phi"
$wb.Worksheets.Item("r.square").Name = "r_square"
$wb.Worksheets.Item("accum.r.square").Name = "accum_r_square"

# --- 2. Flip sign of Component 1 / Component 2 columns on A_hat (rows 2-61, cols B:C) ---
$wsA = $wb.Worksheets.Item("A_hat")
for ($r = 2; $r -le 61; $r++) {
    $wsA.Cells.Item($r, 2).Value2 = -1 * $wsA.Cells.Item($r, 2).Value2
    $wsA.Cells.Item($r, 3).Value2 = -1 * $wsA.Cells.Item($r, 3).Value2
}

# --- 3. Flip sign of Component 1 / Component 2 columns on Phi_hat (rows 2-102, cols A:B) ---
$wsPhi = $wb.Worksheets.Item("Phi_hat")
for ($r = 2; $r -le 102; $r++) {
    $wsPhi.Cells.Item($r, 1).Value2 = -1 * $wsPhi.Cells.Item($r, 1).Value2
    $wsPhi.Cells.Item($r, 2).Value2 = -1 * $wsPhi.Cells.Item($r, 2).Value2
}
